$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "RAD123456"
$ws.Range("F6").Value = "OP12345"
$ws.Range("C8").Value = -692739
$ws.Range("C9").Value = "John Doe"
$ws.Range("C11").Value = "Cargo X"
$ws.Range("C12").Value = "Dependencia Y"
$ws.Range("C13").Value = "CENCO123"
$ws.Range("C14").Value = 100000
$ws.Range("B18").Value = "Descripción de la requisición"
$ws.Range("G30").Value = "'0987654321"
$ws.Range("B43").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B34").Value = "John Doe"
$ws.Range("E34").Value = "Jane Smith"
$ws.Range("B37").Value = "John Doe"
$ws.Range("E37").Value = "Jane Smith"
